# feat: add login component
#
# Mirrors the "Raport zaangazowania" engagement-tracking sheet being
# refreshed with freshly-committed files for the new login component
# (and a CSS file for the pre-existing register component), plus the
# corresponding bump to the existing files' line counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Bump existing "Linie" (lines) counters for files that grew ----
$ws.Range("P37").Value = 3   # was 1
$ws.Range("P40").Value = 7   # was 6
$ws.Range("P42").Value = 6   # was 4

# O42 / O43 lose their (no-op) fill-flag style left over from the prior
# edit session - clear formatting back to the sheet default.
$ws.Range("O42").ClearFormats()
$ws.Range("O43").ClearFormats()

# --- 2. New rows of commits: login component files + register css ----
# Row 50's "N" column had no cell at all yet - borrow the date format from
# the row above so it lines up with N46:N49 instead of landing unformatted.
$ws.Range("N49").Copy()
$ws.Range("N50").PasteSpecial(-4122)  # xlPasteFormats

# Dates first (cheap, order doesn't matter for these).
foreach ($r in 46..50) {
    $ws.Cells.Item($r, 14).Value = 45896   # column N - date
}

# File names - written in first-seen order so new shared-string ids land
# the same way they did in the source commit (login.component.ts,
# login.component.css, login.component.html, register.component.css).
$ws.Cells.Item(46, 15).Value = "login.component.ts"
$ws.Cells.Item(48, 15).Value = "login.component.css"
$ws.Cells.Item(47, 15).Value = "login.component.html"
$ws.Cells.Item(49, 15).Value = "register.component.html"
$ws.Cells.Item(50, 15).Value = "register.component.css"

# Lines-changed counts for each new file.
$ws.Cells.Item(46, 16).Value = 50
$ws.Cells.Item(47, 16).Value = 19
$ws.Cells.Item(48, 16).Value = 46
$ws.Cells.Item(49, 16).Value = 26
$ws.Cells.Item(50, 16).Value = 48

# Match the font used by the other "new file" entries just above (O44:O45)
$newFileRange = $ws.Range("O46:O50")
$newFileRange.Font.Name = "Calibri"
$newFileRange.Font.Size = 11

# --- 3. Refresh the view state (scroll position / active selection) ---
$ws.Range("L53").Select()

$wb.Save()
